$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15
$ws.Range("H15").Value = 1789.2439
$ws.Range("I15").Value = 1789.2439
$ws.Range("K15").Value = 5367.7317
$ws.Range("M15").Value = -5198.7317

# ALC row 68
$ws.Range("H68").Value = 48500.0
$ws.Range("J68").Value = 48500.0
$ws.Range("L68").Value = 48500.0
$ws.Range("N68").Value = -49998.0

# ALC row 71
$ws.Range("H71").Value = 48500.0
$ws.Range("J71").Value = 48500.0
$ws.Range("L71").Value = 145500.0
$ws.Range("N71").Value = -152988.0

# ALC row 116
$ws.Range("H116").Value = 4015.3333
$ws.Range("I116").Value = 3666.7273
$ws.Range("J116").Value = 7850.0
$ws.Range("K116").Value = 3666.7273
$ws.Range("L116").Value = 7850.0
$ws.Range("M116").Value = -224.7273
$ws.Range("N116").Value = -14734.0

# ALC row 118
$ws.Range("H118").Value = 851.1667
$ws.Range("I118").Value = 175.75
$ws.Range("J118").Value = 2202.0
$ws.Range("K118").Value = 527.25
$ws.Range("L118").Value = 6606.0
$ws.Range("M118").Value = 1129.75
$ws.Range("N118").Value = -9920.0

# ALC row 137
$ws.Range("H137").Value = 10824.884
$ws.Range("I137").Value = 3446.9167
$ws.Range("J137").Value = 20144.422
$ws.Range("K137").Value = 10340.7501
$ws.Range("L137").Value = 60433.266
$ws.Range("M137").Value = -7790.750100000001
$ws.Range("N137").Value = -65533.266

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 14696.889
$ws.Range("I32").Value = 6625.469
$ws.Range("K32").Value = 6625.469
$ws.Range("M32").Value = -6338.469

# ARM row 97
$ws.Range("H97").Value = 1644.9375
$ws.Range("I97").Value = 1409.1538
$ws.Range("K97").Value = 1409.1538
$ws.Range("M97").Value = -913.1538

# ARM row 128
$ws.Range("H128").Value = 69999.5
$ws.Range("J128").Value = 69999.5
$ws.Range("L128").Value = 69999.5
$ws.Range("N128").Value = -79959.5

# ARM row 132
$ws.Range("H132").Value = 1287447.9
$ws.Range("I132").Value = 2277.9138
$ws.Range("J132").Value = 5014440.5
$ws.Range("K132").Value = 6833.741399999999
$ws.Range("L132").Value = 15043321.5
$ws.Range("M132").Value = -4303.741399999999
$ws.Range("N132").Value = -15048381.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86
$ws.Range("H86").Value = 4137.846
$ws.Range("I86").Value = 3981.0908
$ws.Range("K86").Value = 3981.0908
$ws.Range("M86").Value = -2858.0908

# BSM row 89
$ws.Range("H89").Value = 4137.846
$ws.Range("I89").Value = 3981.0908
$ws.Range("K89").Value = 19905.454
$ws.Range("M89").Value = -14289.454

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 10309.792
$ws.Range("I31").Value = 8594.5
$ws.Range("J31").Value = 10759.705
$ws.Range("K31").Value = 8594.5
$ws.Range("L31").Value = 10759.705
$ws.Range("M31").Value = -8299.5
$ws.Range("N31").Value = -11349.705

# CRP row 34
$ws.Range("H34").Value = 10309.792
$ws.Range("I34").Value = 8594.5
$ws.Range("J34").Value = 10759.705
$ws.Range("K34").Value = 8594.5
$ws.Range("L34").Value = 10759.705
$ws.Range("M34").Value = -8392.5
$ws.Range("N34").Value = -11163.705

# CRP row 74
$ws.Range("H74").Value = 37500.0
$ws.Range("J74").Value = 37500.0
$ws.Range("L74").Value = 37500.0
$ws.Range("N74").Value = -39248.0

# CRP row 77
$ws.Range("H77").Value = 37500.0
$ws.Range("J77").Value = 37500.0
$ws.Range("L77").Value = 112500.0
$ws.Range("N77").Value = -121236.0

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5
$ws.Range("H5").Value = 1628157.1
$ws.Range("I5").Value = 1115.0
$ws.Range("J5").Value = 2441678.0
$ws.Range("K5").Value = 3345.0
$ws.Range("L5").Value = 7325034.0
$ws.Range("M5").Value = -3233.0
$ws.Range("N5").Value = -7325258.0

# CUL row 107
$ws.Range("H107").Value = 2084158.5
$ws.Range("I107").Value = 422.25
$ws.Range("J107").Value = 2841880.8
$ws.Range("K107").Value = 1266.75
$ws.Range("L107").Value = 8525642.399999999
$ws.Range("M107").Value = 653.25
$ws.Range("N107").Value = -8529482.399999999

# CUL row 119
$ws.Range("H119").Value = 8999.5
$ws.Range("J119").Value = 0.0
$ws.Range("L119").Value = 0.0
$ws.Range("N119").ClearContents()

# CUL row 129
$ws.Range("H129").Value = 1123.3928
$ws.Range("I129").Value = 756.6923
$ws.Range("K129").Value = 2270.0769
$ws.Range("M129").Value = 2729.9231

# CUL row 131
$ws.Range("H131").Value = 1461.17
$ws.Range("I131").Value = 1069.6
$ws.Range("J131").Value = 1481.7789
$ws.Range("K131").Value = 3208.8
$ws.Range("L131").Value = 4445.3367
$ws.Range("M131").Value = 1831.2
$ws.Range("N131").Value = -14525.3367

# CUL row 132
$ws.Range("H132").Value = 1506.1428
$ws.Range("I132").Value = 1290.8572
$ws.Range("J132").Value = 1721.4286
$ws.Range("K132").Value = 11617.7148
$ws.Range("L132").Value = 15492.8574
$ws.Range("M132").Value = -9087.7148
$ws.Range("N132").Value = -20552.8574

# CUL row 135
$ws.Range("H135").Value = 1628157.1
$ws.Range("I135").Value = 1115.0
$ws.Range("J135").Value = 2441678.0
$ws.Range("K135").Value = 10035.0
$ws.Range("L135").Value = 21975102.0
$ws.Range("M135").Value = -7500.0
$ws.Range("N135").Value = -21980172.0

# CUL row 136
$ws.Range("H136").Value = 2574.0
$ws.Range("I136").Value = 2574.0
$ws.Range("J136").Value = 0.0
$ws.Range("K136").Value = 7722.0
$ws.Range("L136").Value = 0.0
$ws.Range("M136").Value = -2622.0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# GSM row 97
$ws.Range("H97").Value = 2039.2424
$ws.Range("I97").Value = 1074.1305
$ws.Range("J97").Value = 4259.0
$ws.Range("K97").Value = 1074.1305
$ws.Range("L97").Value = 4259.0
$ws.Range("M97").Value = -578.1305
$ws.Range("N97").Value = -5251.0

# GSM row 113
$ws.Range("H113").Value = 1332.2106
$ws.Range("I113").Value = 1122.3572
$ws.Range("K113").Value = 1122.3572
$ws.Range("M113").Value = 1047.6428

# GSM row 122
$ws.Range("H122").Value = 2407.0
$ws.Range("I122").Value = 2009.25
$ws.Range("J122").Value = 3003.625
$ws.Range("K122").Value = 6027.75
$ws.Range("L122").Value = 9010.875
$ws.Range("M122").Value = -3577.75
$ws.Range("N122").Value = -13910.875

# GSM row 126
$ws.Range("H126").Value = 3884.4736
$ws.Range("I126").Value = 3162.6296
$ws.Range("J126").Value = 5656.273
$ws.Range("K126").Value = 9487.8888
$ws.Range("L126").Value = 16968.819
$ws.Range("M126").Value = -7017.888800000001
$ws.Range("N126").Value = -21908.819

$ws = $wb.Worksheets.Item("LTW")
# LTW row 74
$ws.Range("H74").Value = 47250.0
$ws.Range("I74").Value = 47250.0
$ws.Range("K74").Value = 47250.0
$ws.Range("M74").Value = -46252.0

# LTW row 77
$ws.Range("H77").Value = 47250.0
$ws.Range("I77").Value = 47250.0
$ws.Range("K77").Value = 141750.0
$ws.Range("M77").Value = -136758.0

# LTW row 131
$ws.Range("H131").Value = 0.0
$ws.Range("I131").Value = 0.0
$ws.Range("K131").Value = 0.0
$ws.Range("M131").ClearContents()

# LTW row 135
$ws.Range("H135").Value = 95833.336
$ws.Range("J135").Value = 95833.336
$ws.Range("L135").Value = 95833.336
$ws.Range("N135").Value = -105973.336

$ws = $wb.Worksheets.Item("WVR")
# WVR row 2
$ws.Range("H2").Value = 771602.56
$ws.Range("I2").Value = 771602.56
$ws.Range("K2").Value = 771602.56
$ws.Range("M2").Value = -771490.56

Write-Host "Applied all Leve profit updates."
